$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111363024
$ws.Range("B2").Value = 77268
$ws.Range("E2").Value = 228912
$ws.Range("F2").Value = "Mörk kolflarnlav"
$ws.Range("G2").Value = "Carbonicola myrmecina"
$ws.Range("H2").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q2").Value = 593291.0260186956
$ws.Range("R2").Value = 6987171.95495991

# Row 4
$ws.Range("A4").Value = 111363026
$ws.Range("B4").Value = 90854
$ws.Range("E4").Value = 2079
$ws.Range("F4").Value = "Nordtagging"
$ws.Range("G4").Value = "Odonticium romellii"
$ws.Range("H4").Value = "(S.Lundell) Parmasto"
$ws.Range("Q4").Value = 593292.3890792141
$ws.Range("R4").Value = 6987203.815111163
$ws.Range("S4").Value = 10

# Row 5
$ws.Range("A5").Value = 111363028
$ws.Range("B5").Value = 77186
$ws.Range("E5").Value = 353
$ws.Range("F5").Value = "Dvärgbägarlav"
$ws.Range("G5").Value = "Cladonia parasitica"
$ws.Range("H5").Value = "(Hoffm.) Hoffm."
$ws.Range("Q5").Value = 593324.0129203054
$ws.Range("R5").Value = 6987101.07452714

# Row 6
$ws.Range("A6").Value = 111363030
$ws.Range("B6").Value = 77268
$ws.Range("E6").Value = 228912
$ws.Range("F6").Value = "Mörk kolflarnlav"
$ws.Range("G6").Value = "Carbonicola myrmecina"
$ws.Range("H6").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q6").Value = 593355.1995546351
$ws.Range("R6").Value = 6987156.520171621
$ws.Range("S6").Value = 25

# Row 7
$ws.Range("A7").Value = 111363022
$ws.Range("B7").Value = 77186
$ws.Range("E7").Value = 353
$ws.Range("F7").Value = "Dvärgbägarlav"
$ws.Range("G7").Value = "Cladonia parasitica"
$ws.Range("H7").Value = "(Hoffm.) Hoffm."
$ws.Range("Q7").Value = 593324.9051589288
$ws.Range("R7").Value = 6987181.108611984

# Row 8
$ws.Range("A8").Value = 111363020
$ws.Range("B8").Value = 78107
$ws.Range("E8").Value = 6453
$ws.Range("F8").Value = "Vedskivlav"
$ws.Range("G8").Value = "Hertelidea botryosa"
$ws.Range("H8").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q8").Value = 593324.7367794912
$ws.Range("R8").Value = 6987171.102828567

# Row 9
$ws.Range("A9").Value = 111363021
$ws.Range("B9").Value = 89330
$ws.Range("E9").Value = 3242
$ws.Range("F9").Value = "Vitplätt"
$ws.Range("G9").Value = "Chaetodermella luna"
$ws.Range("H9").Value = "(Romell ex D.P.Rogers & H.S.Jacks.) Rauschert"
$ws.Range("Q9").Value = 593278.356042281
$ws.Range("R9").Value = 6987153.408284122

# Row 10
$ws.Range("A10").Value = 111363031
$ws.Range("B10").Value = 76918
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 6437
$ws.Range("F10").Value = "Blanksvart spiklav"
$ws.Range("G10").Value = "Calicium denigratum"
$ws.Range("H10").Value = "(Vain.) Tibell"
$ws.Range("Q10").Value = 593417.4633552339
$ws.Range("R10").Value = 6986985.556671137

# Row 11
$ws.Range("A11").Value = 111363025
$ws.Range("B11").Value = 89646
$ws.Range("D11").Value = "VU"
$ws.Range("E11").Value = 65
$ws.Range("F11").Value = "Fläckporing"
$ws.Range("G11").Value = "Anthoporia albobrunnea"
$ws.Range("H11").Value = "(Romell) Karasiński & Niemelä"
$ws.Range("Q11").Value = 593292.3890792141
$ws.Range("R11").Value = 6987203.815111163

# Row 12
$ws.Range("A12").Value = 111363029
$ws.Range("B12").Value = 76918
$ws.Range("E12").Value = 6437
$ws.Range("F12").Value = "Blanksvart spiklav"
$ws.Range("G12").Value = "Calicium denigratum"
$ws.Range("H12").Value = "(Vain.) Tibell"
$ws.Range("Q12").Value = 593291.0260186956
$ws.Range("R12").Value = 6987171.95495991
